# Rewrites the ensemble-model metrics sheet: every model row (2-26) now
# carries the same, newly-computed metric values (columns B:Q), and the
# model-name labels in column A are reassigned/reshuffled per the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values shared by every data row (B:Q), in column order.
$values = @(0.9999338618126242, 0.9991272753483447, 0.9999748392634773, 0.9999698392805606, 0.9999724124880138, 0.00006173709478443979, 0.000814650154135481, 0.0000255670867766081, 0.00003428362379971795, 0.00002992535528816303, 0.0004885546030702877, 0.007857295640641236, 1.000317463299404, 0.008191796982907578, 77.38525119128605, 112.7326501124639)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

for ($r = 2; $r -le 26; $r++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}

# New model-name labels (column A), keyed by row number.
$names = @{
    2  = "model_10_9_0"
    3  = "model_10_9_22"
    4  = "model_10_9_21"
    5  = "model_10_9_20"
    6  = "model_10_9_19"
    7  = "model_10_9_18"
    8  = "model_10_9_17"
    9  = "model_10_9_16"
    10 = "model_10_9_15"
    11 = "model_10_9_14"
    12 = "model_10_9_13"
    13 = "model_10_9_23"
    14 = "model_10_9_12"
    15 = "model_10_9_10"
    16 = "model_10_9_9"
    17 = "model_10_9_8"
    18 = "model_10_9_7"
    19 = "model_10_9_6"
    20 = "model_10_9_5"
    21 = "model_10_9_4"
    22 = "model_10_9_3"
    23 = "model_10_9_2"
    24 = "model_10_9_1"
    25 = "model_10_9_11"
    26 = "model_10_9_24"
}

foreach ($r in $names.Keys) {
    $ws.Range("A" + $r).Value = $names[$r]
}
